$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Add the new J column values for rows 3 and 4 (mirroring H/I data)
$ws.Range("J3").Value = 7
$ws.Range("J4").Value = 28

# Update the active selection to L5 (was L9)
$ws.Range("L5").Select()
